$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update selection (active cell) on the sheet view
$ws.Range("W22").Select()

# Row 16 - numeric values (min)
$ws.Range("T16").Value = 56
$ws.Range("U16").Value = 55
$ws.Range("V16").Value = 55
$ws.Range("W16").Value = 55

# Row 17 - text values (max, pulled from shared strings list: 77, 77, 79)
$ws.Range("U17").Value = "77"
$ws.Range("V17").Value = "77"
$ws.Range("W17").Value = "79"

# Row 18 - % in 50s
$ws.Range("T18").Value = 0.13
$ws.Range("U18").Value = 0.19
$ws.Range("V18").Value = 0.15
$ws.Range("W18").Value = 0.19

# Row 19 - % in 60s
$ws.Range("T19").Value = 0.78
$ws.Range("U19").Value = 0.71
$ws.Range("V19").Value = 0.65
$ws.Range("W19").Value = 0.63

# Row 20 - % in 70s
$ws.Range("T20").Value = 0.09
$ws.Range("U20").Value = 0.1
$ws.Range("V20").Value = 0.21
$ws.Range("W20").Value = 0.19

# Row 21 - % in 80s
$ws.Range("T21").Value = 0
$ws.Range("U21").Value = 0
$ws.Range("V21").Value = 0
$ws.Range("W21").Value = 0

# Row 22 - Objective / DONE markers
$ws.Range("T22").Value = "DONE"
$ws.Range("U22").Value = "DONE"
$ws.Range("V22").Value = "DONE"
$ws.Range("W22").Value = "DONE"
